$d = $word.ActiveDocument

# Locate the part of the sentence that must stay ("...cada DNI") and the
# part that must be dropped (" y la evaluación...Parte 1.").
$keepText = "Implementó el conteo de frecuencia de los dígitos, la suma total de los dígitos de cada DNI"
$dropText = " y la evaluación de condiciones lógicas basadas en las expresiones trabajadas en la Parte 1."

$findRange = $d.Content
$found = $findRange.Find.Execute($keepText)
if (-not $found) {
    throw "Could not locate the target sentence."
}
$keepEnd = $findRange.End

# Isolate the trailing portion we are about to remove into its own run
# (a harmless format round-trip forces a run split) without touching the
# neighbouring run(s) elsewhere in the paragraph.
$tail = $d.Range($keepEnd, $keepEnd + $dropText.Length)
if ($tail.Text -ne $dropText) {
    throw "Unexpected tail text: $($tail.Text)"
}
$tail.Font.Size = 99
$tail.Font.Size = 12

# Now delete the isolated trailing run's text.
$tail2 = $d.Range($keepEnd, $keepEnd + $dropText.Length)
$tail2.Text = ""

# Insert the closing period as its own run right after the kept text.
$insertionPoint = $d.Range($keepEnd, $keepEnd)
$insertionPoint.InsertAfter(".")

# Match the surrounding paragraph's font size (w:sz 24 half-points = 12pt).
$period = $d.Range($keepEnd, $keepEnd + 1)
$period.Font.Size = 99
$period.Font.Size = 12

$d.Save()
